# Update Participant Questionnaire Excel
# Insert 9 new participant rows (P01-P09) above the existing P10-P18 rows on
# "Tabelle1", populate the new rows with their data, resize column G, and
# update the sheet view (zoom + selection) to match the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert 9 blank rows before the current row 2 - this pushes the existing
# participants (currently P10..P18 on rows 2-10) down to rows 11-19.
$ws.Rows("2:10").Insert()

# Fill in the new participant rows. Values are entered from the bottom
# (P09) up to the top (P01) so that the new strings land in the shared
# string table in that same order.
$ws.Range("A10").Value = "P09"

$ws.Range("A9").Value = "P08"
$ws.Range("B9").Value = 26
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1

$ws.Range("A8").Value = "P07"
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0

$ws.Range("A7").Value = "P06"
$ws.Range("A6").Value = "P05"
$ws.Range("A5").Value = "P04"
$ws.Range("A4").Value = "P03"
$ws.Range("A3").Value = "P02"
$ws.Range("A2").Value = "P01"

# Widen column G (it is no longer sized by "best fit").
$ws.Columns("G:G").ColumnWidth = 5.75

# Update the view: zoom to 150% and select C13:C14.
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("C13:C14").Select()
